$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells (stored as real numbers, type "n") ---
$ws.Range("A3").Value = 131274492
$ws.Range("B3").Value = 57725
$ws.Range("E3").Value = 102621
$ws.Range("Q3").Value = 697494
$ws.Range("R3").Value = 6640927
$ws.Range("S3").Value = 4

# --- Plain text cells (not numeric-looking, safe as-is) ---
$ws.Range("D3").Value = "LC"
$ws.Range("F3").Value = "Sparvuggla"
$ws.Range("G3").Value = "Glaucidium passerinum"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("P3").Value = "Ladängssjöns naturreservat, Upl"
$ws.Range("T3").Value = "Stockholm"
$ws.Range("U3").Value = "Norrtälje"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Estuna"
$ws.Range("Z3").Value = "16:40"
$ws.Range("AB3").Value = "16:40"
$ws.Range("AC3").Value = "Upptäcktes från bilen. Tyst."
$ws.Range("AW3").Value = "Peter Border"
$ws.Range("AX3").Value = "Peter Border, Olle Rådfeldt"

# --- Text cells whose content looks numeric/date, force text format first
#     so the engine keeps them literal instead of coercing to number/date ---
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2026-02-23"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2026-02-23"

# --- Boolean cells ---
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# --- Present-but-empty text cells (mirrors the existing blank cells in
#     row 2, e.g. AT2/AY2, which are empty text/string cells). A lone
#     quote is Excel's "text, no contents" entry -- it yields a real,
#     present cell of text type with an empty value, same as the source
#     file, instead of being dropped like an ordinary empty assignment. ---
$ws.Range("K3").Value = "'"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("AT3").Value = "'"
$ws.Range("AY3").Value = "'"
